$d = $word.ActiveDocument

# --- Edit 1: remove the stray _GoBack bookmark that sits after "which" -----
# A genuine Find/Replace across the span that contains the bookmark rewrites
# that part of the paragraph and drops the (now orphaned) bookmark markers.
$d.Content.Find.Execute("which enrolled", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "which enrolled", 2)

# --- Edit 2: rewrite the "2. The system shall ..." requirement sentence ----
$oldReq = "2. The system shall display details or information of the course, including course name, course ID, course seat, credit, semester, year and teacher names."
$newReq = "2. The system shall display details or information of the course, including Course name, Course ID, Course seat, Credit, Semester, Year and Teacher names."

$d.Content.Find.Execute($oldReq, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newReq, 2)

# Re-insert the "_GoBack" bookmark right after "2. The system shall " (and
# before "display ..."), which is where it now belongs per the new text.
$r = $d.Content
$r.Find.Execute("2. The system shall ", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$splitPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $splitPoint)
